$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$nc = @(19330051920042, 19330051920049, 19330051920056, 19330051920068, 19330051920081, 19330051920037)
$paterno = @("ARIZMENDI", "CRUZ", "GONZALEZ", "MARTINEZ", "VAZQUEZ", "VALENCIA")
$materno = @("NUÑEZ", "GALVEZ", "DE LOS SANTOS", "PACHECO", "ROMERO", "GARCIA")
$nombres = @("KARLA IRAN", "LUIS ARIEL", "MARTIN", "EMMANUEL", "MONSERRAT", "MAURICIO")
$nombreLargo = "CIENCIA, TECNOLOGÍA, SOCIEDAD Y VALORES"
$grupo = @("5AEV", "5AEV", "5AEV", "5AEV", "5AEV", "5AEM")
$reprobadas = @(7, 7, 7, 7, 7, 7)

for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $nc[$i]
}
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $paterno[$i]
}
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $materno[$i]
}
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $nombres[$i]
}
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $nombreLargo
}
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $grupo[$i]
}
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $reprobadas[$i]
}
